$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = '61.972.76'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '3.408.65'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  -0.06%  '
Set-TextValue "D5" '409.35'
$ws.Range("E5").Value = '  +0.58%  '
Set-TextValue "D6" '129.49'
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("E7").Value = '  +6.37%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +5.64%  '
Set-TextValue "D10" '0.142'
$ws.Range("E10").Value = '  +2.67%  '
$ws.Range("E11").Value = '  +2.03%  '
$ws.Range("E12").Value = '  +42.37%  '
$ws.Range("E13").Value = '  +10.72%  '
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = '3.953.33'
$ws.Range("E15").Value = '  -0.63%  '
Set-TextValue "D16" '21.20'
$ws.Range("D17").Value = '3.407.54'
$ws.Range("E17").Value = '  -0.65%  '
Set-TextValue "D18" '12.55'
$ws.Range("E18").Value = '  +9.16%  '
$ws.Range("E19").Value = '  +6.99%  '
$ws.Range("D20").Value = '61.980.98'
$ws.Range("E20").Value = '  -1.15%  '
Set-TextValue "D21" '451.20'
$ws.Range("E21").Value = '  +43.00%  '
Set-TextValue "D22" '92.03'
$ws.Range("E22").Value = '  +8.97%  '
$ws.Range("E23").Value = '  +1.17%  '
Set-TextValue "D24" '13.23'
$ws.Range("E24").Value = '  +3.23%  '
Set-TextValue "D25" '3.30'
$ws.Range("E25").Value = '  +3.84%  '
Set-TextValue "D26" '9.35'
$ws.Range("E26").Value = '  +14.67%  '
$ws.Range("E27").Value = '  +11.63%  '
Set-TextValue "D28" '4.78'
$ws.Range("E28").Value = '  +0.30%  '
Set-TextValue "D29" '7.74'
$ws.Range("E29").Value = '  -1.12%  '
$ws.Range("E30").Value = '  +0.82%  '
Set-TextValue "D31" '12.00'
$ws.Range("E31").Value = '  +5.35%  '
$ws.Range("E32").Value = '  -1.71%  '
$ws.Range("E33").Value = '  +0.05%  '
Set-TextValue "D34" '42.87'
$ws.Range("E34").Value = '  -3.62%  '
$ws.Range("E35").Value = '  -0.01%  '
Set-TextValue "D36" '0.0506'
$ws.Range("E36").Value = '  +4.45%  '
Set-TextValue "D37" '53.70'
$ws.Range("E37").Value = '  +3.54%  '
Set-TextValue "D38" '0.999'
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("E39").Value = '  +1.81%  '
$ws.Range("E40").Value = '  +7.75%  '
$ws.Range("E41").Value = '  -0.53%  '
Set-TextValue "D42" '0.321'
$ws.Range("E42").Value = '  -0.86%  '
Set-TextValue "D43" '143.11'
$ws.Range("E43").Value = '  +0.26%  '
Set-TextValue "D44" '4.28'
$ws.Range("E44").Value = '  +9.50%  '
$ws.Range("E45").Value = '  +16.13%  '
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("E47").Value = '  -1.26%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D48" '0.148'
$ws.Range("E48").Value = '  +22.72%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D49" '22.34'
$ws.Range("E49").Value = '  +4.91%  '
$ws.Range("E50").Value = '  +9.02%  '
$ws.Range("B51").Value = 'Fetch.AI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D51" '1.93'
$ws.Range("E51").Value = '  +14.89%  '
